$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff"
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date / Datetime updates
$overview.Range("D2").Value = "2016-03-24 09:39:15"
$dede.Range("E2").Value = "2016-03-24 09:39:15"
$zhcn.Range("E2").Value = "2016-03-24 09:39:11"
